# [IMP] Transdoo upgrade (PO state) & test
#
# The translation export sheet ("odoo_default_tnl") gained one new
# msgid/msgstr pair: "Total Due" / "Totale dovuto". It needs to be
# inserted in its alphabetically-sorted position among the "Total ..."
# entries, i.e. right after "Total debit vat" / "Totale IVA a debito"
# (old row 1447) and before "Total in Company Currency" (old row 1448).
# Inserting the row there pushes every following row down by one, so the
# sheet grows from 1579 to 1580 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("odoo_default_tnl")

# Insert a new blank row at position 1448; everything from the old row
# 1448 onward shifts down to 1449+.
$ws.Rows.Item(1448).Insert()

# Fill in the new translation pair (column A = status, B = msgid,
# C = msgstr; this row only uses B/C like its neighbours).
$ws.Cells.Item(1448, 2).Value = "Total Due"
$ws.Cells.Item(1448, 3).Value = "Totale dovuto"

# Best-effort: reflect the post-edit selection/viewport (cosmetic only,
# matches the author having scrolled to and selected the full updated
# range after the edit).
$lastRow = $ws.UsedRange.Rows.Count
$ws.Range("A1:C" + $lastRow).Select() | Out-Null
